$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("source")

$ws.Range("A1").Value = "Generated By: \\pfs1w\C:\Users\kheal579\Documents\01_integral-private\integral-private\vignettes\C:\Users\kheal579\AppData\Local\Temp\Rtmp2nWN2H\callr-scr-99a4716e1e26"
$ws.Range("A2").Value = "Created By: Eben Pendleton on 2022-08-24"
